# update R values and dates
#
# The "R0" sheet holds a small time-series of R0 estimates keyed by date.
# Two stale data points are removed, the most recent pre-existing estimate
# is corrected, and a new, more recent estimate is appended - replacing the
# previous last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("R0")

# Drop the 1/02/2020 (43862) row - no longer part of the series.
$ws.Rows(3).Delete() | Out-Null

# Drop the 12/03/2020 (43902) row - superseded.
$ws.Rows(4).Delete() | Out-Null

# The 18/03/2020 estimate is revised from 0.99 to 1.
$ws.Range("C5").Value = 1

# Replace the old trailing 2/04/2020 (43923/0.972) row in place with the
# new 7/04/2020 (43928/0.85) estimate, then drop the now-duplicated old
# trailing row that got shifted down.
$ws.Range("B6").Value = 43928
$ws.Range("C6").Value = 0.85
$ws.Rows(7).Delete() | Out-Null

# Leave a few formatted-but-empty rows at the very bottom of the sheet,
# mirroring the trailing blank rows left behind in the sheet.
$ws.Rows(1048573).RowHeight = 12.8
$ws.Rows(1048574).RowHeight = 12.8
$ws.Rows(1048575).RowHeight = 12.8
$ws.Rows(1048576).RowHeight = 12.8

# Cursor ends up back near the top of the refreshed data.
$ws.Range("A6").Select() | Out-Null
